# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Ultros_Profits workbook (FFXIV leve-profit sheets)
# per the authoritative diff (per-sheet, per-cell <v> changes).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1237.1111
$ws.Range("J33").Value = 1197.5
$ws.Range("L33").Value = 1197.5
$ws.Range("N33").Value = -1655.5
$ws.Range("H37").Value = 6000
$ws.Range("I37").Value = 6000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 18000
$ws.Range("N37").ClearContents()
$ws.Range("M37").Value = -17874
$ws.Range("H40").Value = 5000.8076
$ws.Range("I40").Value = 3100
$ws.Range("J40").Value = 5159.2085
$ws.Range("K40").Value = 3100
$ws.Range("L40").Value = 5159.2085
$ws.Range("M40").Value = -2925
$ws.Range("N40").Value = -5509.2085
$ws.Range("H101").Value = 560.8570999999999
$ws.Range("I101").Value = 571
$ws.Range("K101").Value = 1713
$ws.Range("M101").Value = -91
$ws.Range("H112").Value = 1733.8572
$ws.Range("J112").Value = 1740.1708
$ws.Range("L112").Value = 5220.512400000001
$ws.Range("N112").Value = -7436.512400000001
$ws.Range("H125").Value = 15540.286
$ws.Range("I125").Value = 1582.6666
$ws.Range("K125").Value = 14243.9994
$ws.Range("M125").Value = -11783.9994
$ws.Range("H138").Value = 2331.2903
$ws.Range("I138").Value = 1312.3636
$ws.Range("J138").Value = 4822
$ws.Range("K138").Value = 3937.0908
$ws.Range("L138").Value = 14466
$ws.Range("M138").Value = 1202.9092
$ws.Range("N138").Value = -24746
$ws.Range("H141").Value = 6738.4287
$ws.Range("I141").Value = 7048.1665
$ws.Range("K141").Value = 21144.4995
$ws.Range("M141").Value = -15964.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12722
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 10723.712
$ws.Range("I32").Value = 8317.057000000001
$ws.Range("K32").Value = 8317.057000000001
$ws.Range("M32").Value = -8030.057000000001
$ws.Range("H61").Value = 4336.4614
$ws.Range("I61").Value = 3068.4614
$ws.Range("K61").Value = 3068.4614
$ws.Range("M61").Value = -2856.4614
$ws.Range("H116").Value = 12722
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -9588
$ws.Range("H136").Value = 4336.4614
$ws.Range("I136").Value = 3068.4614
$ws.Range("K136").Value = 9205.3842
$ws.Range("M136").Value = -6655.3842
$ws.Range("H139").Value = 80715
$ws.Range("J139").Value = 80715
$ws.Range("L139").Value = 80715
$ws.Range("N139").Value = -90995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12722
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H20").Value = 3012.9
$ws.Range("I20").Value = 3071
$ws.Range("K20").Value = 3071
$ws.Range("M20").Value = -2824
$ws.Range("H22").Value = 3824.8572
$ws.Range("I22").Value = 4754.4
$ws.Range("K22").Value = 4754.4
$ws.Range("M22").Value = -4581.4
$ws.Range("H86").Value = 5507
$ws.Range("I86").Value = 1389.8572
$ws.Range("J86").Value = 11271
$ws.Range("K86").Value = 1389.8572
$ws.Range("L86").Value = 11271
$ws.Range("M86").Value = -266.8571999999999
$ws.Range("N86").Value = -13517
$ws.Range("H89").Value = 5507
$ws.Range("I89").Value = 1389.8572
$ws.Range("J89").Value = 11271
$ws.Range("K89").Value = 6949.286
$ws.Range("L89").Value = 56355
$ws.Range("M89").Value = -1333.286
$ws.Range("N89").Value = -67587
$ws.Range("H130").Value = 49985.938
$ws.Range("J130").Value = 49985.938
$ws.Range("L130").Value = 49985.938
$ws.Range("N130").Value = -60025.938
$ws.Range("H134").Value = 5254
$ws.Range("I134").Value = 3692.7273
$ws.Range("K134").Value = 11078.1819
$ws.Range("M134").Value = -8543.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 285.75
$ws.Range("I39").Value = 285.75
$ws.Range("K39").Value = 285.75
$ws.Range("M39").Value = 105.25
$ws.Range("H49").Value = 285.75
$ws.Range("I49").Value = 285.75
$ws.Range("K49").Value = 285.75
$ws.Range("M49").Value = -103.75
$ws.Range("H59").Value = 31933.572
$ws.Range("J59").Value = 34755.832
$ws.Range("L59").Value = 34755.832
$ws.Range("N59").Value = -37045.832
$ws.Range("H122").Value = 41669404
$ws.Range("I122").Value = 71429860
$ws.Range("K122").Value = 214289580
$ws.Range("M122").Value = -214287130
$ws.Range("H132").Value = 4416.5
$ws.Range("I132").Value = 2972
$ws.Range("K132").Value = 8916
$ws.Range("M132").Value = -6386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3543.625
$ws.Range("I3").Value = 3058.1667
$ws.Range("K3").Value = 9174.500100000001
$ws.Range("M3").Value = -9062.500100000001
$ws.Range("H40").Value = 56.272728
$ws.Range("I40").Value = 49.5
$ws.Range("K40").Value = 198
$ws.Range("M40").Value = -129
$ws.Range("H131").Value = 5974.625
$ws.Range("J131").Value = 7502.1113
$ws.Range("L131").Value = 22506.3339
$ws.Range("N131").Value = -32586.3339
$ws.Range("H132").Value = 683.3333
$ws.Range("I132").Value = 550
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 4950
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -2420
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 15198.6
$ws.Range("J40").Value = 15198.6
$ws.Range("L40").Value = 15198.6
$ws.Range("N40").Value = -15500.6
$ws.Range("H43").Value = 9380
$ws.Range("H132").Value = 4664.926
$ws.Range("I132").Value = 3788.8823
$ws.Range("J132").Value = 6154.2
$ws.Range("K132").Value = 11366.6469
$ws.Range("L132").Value = 18462.6
$ws.Range("M132").Value = -8836.6469
$ws.Range("N132").Value = -23522.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4204.9414
$ws.Range("I7").Value = 2183.1428
$ws.Range("K7").Value = 2183.1428
$ws.Range("M7").Value = -2071.1428
$ws.Range("H22").Value = 647.2222
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 625
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 625
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1215
$ws.Range("H27").Value = 647.2222
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 625
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 625
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -839
$ws.Range("H61").Value = 2297.0344
$ws.Range("I61").Value = 1668.2
$ws.Range("J61").Value = 2970.7856
$ws.Range("K61").Value = 1668.2
$ws.Range("L61").Value = 2970.7856
$ws.Range("M61").Value = -1466.2
$ws.Range("N61").Value = -3374.7856
$ws.Range("H82").Value = 4320.6
$ws.Range("I82").Value = 2237.5
$ws.Range("J82").Value = 8486.799999999999
$ws.Range("K82").Value = 2237.5
$ws.Range("L82").Value = 8486.799999999999
$ws.Range("M82").Value = -1876.5
$ws.Range("N82").Value = -9208.799999999999
$ws.Range("H85").Value = 4320.6
$ws.Range("I85").Value = 2237.5
$ws.Range("J85").Value = 8486.799999999999
$ws.Range("K85").Value = 2237.5
$ws.Range("L85").Value = 8486.799999999999
$ws.Range("M85").Value = -989.5
$ws.Range("N85").Value = -10982.8
$ws.Range("H100").Value = 131345.67
$ws.Range("I100").Value = 1111111
$ws.Range("K100").Value = 1111111
$ws.Range("M100").Value = -1110570
$ws.Range("H113").Value = 2297.0344
$ws.Range("I113").Value = 1668.2
$ws.Range("J113").Value = 2970.7856
$ws.Range("K113").Value = 1668.2
$ws.Range("L113").Value = 2970.7856
$ws.Range("M113").Value = 501.8
$ws.Range("N113").Value = -7310.7856
$ws.Range("H126").Value = 4204.9414
$ws.Range("I126").Value = 2183.1428
$ws.Range("K126").Value = 6549.428400000001
$ws.Range("M126").Value = -4079.428400000001
$ws.Range("H132").Value = 4112.4
$ws.Range("I132").Value = 3365.1
$ws.Range("K132").Value = 10095.3
$ws.Range("M132").Value = -7565.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2215.587
$ws.Range("I132").Value = 1628.6666
$ws.Range("J132").Value = 4328.5
$ws.Range("K132").Value = 4885.9998
$ws.Range("L132").Value = 12985.5
$ws.Range("M132").Value = -2355.9998
$ws.Range("N132").Value = -18045.5
$ws.Range("H136").Value = 4412.1143
$ws.Range("J136").Value = 5127
$ws.Range("L136").Value = 15381
$ws.Range("N136").Value = -20481
